$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the alias for Augusta Delono (row 4) in the "Alias" column (D)
# and stash the previous alias value off to the side in column J.
$ws.Range("J4").Value = "AUD"
$ws.Range("D4").Value = "ADO"
